$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# "FLD_DocumentRegistry_New_Transmittals" is being (re-)included and marked
# for execution, so its stale "PASS" result is cleared back out (column
# "Results", row 2 -> cell E2).
$ws.Range("E2").ClearContents()

# Drop the stale selection state (previously parked on A2) back to the
# sheet's default A1 position.
$ws.Range("A1").Select()
